# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" rows (16-22) are re-sorted ascending by period.
# Originally the table listed periods descending (2006..1912); the
# updated workbook lists them ascending (1912..2006). Column B-D and G-J
# are identical across these rows, so the practical effect of the sort
# is limited to swapping the "Periodo Mora" (E) and "Valor Mora" (F)
# values between the mirrored row pairs: 16<->22, 17<->21, 18<->20.
# Row 19 (period 2003) already sits in its sorted position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-PeriodoRows($rowA, $rowB) {
    $eA = $ws.Range("E$rowA").Value()
    $fA = $ws.Range("F$rowA").Value()
    $eB = $ws.Range("E$rowB").Value()
    $fB = $ws.Range("F$rowB").Value()

    $ws.Range("E$rowA").Value = $eB
    $ws.Range("F$rowA").Value = $fB
    $ws.Range("E$rowB").Value = $eA
    $ws.Range("F$rowB").Value = $fA
}

Swap-PeriodoRows 16 22
Swap-PeriodoRows 17 21
Swap-PeriodoRows 18 20
# Row 19 unchanged (period 2003 already in its sorted spot).
